$d = $word.ActiveDocument

# --- Send each floating picture backward one step in the z-order stack.
# This mirrors what Word does to relativeHeight (wp:anchor) when the
# drawing stack is touched/resaved: each of the 4 anchors' relativeHeight
# drops by 1024, preserving their relative order.
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    [void]$shp.ZOrder(3)  # msoSendBackward
}

# --- Add the new "Hash tables + Hashing" section at the end of the doc,
# right before the trailing empty paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertion = $lastPara.Range
$insertion.Collapse(1)  # wdCollapseStart
$insertion.InsertBefore("Hash tables + Hashing:`r")

# Format the heading paragraph (text + paragraph mark) with
# sz/szCs 24 (12pt) and a single underline, matching the other section
# headings used throughout the document.
$headingPara = $d.Paragraphs.Item($count)
$headingRange = $headingPara.Range
$headingRange.Font.Size = 12
$headingRange.Font.SizeBi = 12
$headingRange.Font.Underline = 1

# Insert the body paragraph describing hash tables, right before the
# trailing empty paragraph.
$count2 = $d.Paragraphs.Count
$lastPara2 = $d.Paragraphs.Item($count2)
$insertion2 = $lastPara2.Range
$insertion2.Collapse(1)
$bodyText = "Hash tables are data structures that create maps between keys and values. They use a hashing algorithm to calculate the index a certain piece of data, or " + [char]0x2018 + "key" + [char]0x2019 + ", should be placed into in the table. This immensely speeds up the search for any item of data in the table as the " + [char]0x2018 + "key" + [char]0x2019 + " can be put into the hashing algorithm to calculate its index, and then the data at that index in the table can be found immediately. Hashing algorithms are a calculation that takes place on the " + [char]0x2018 + "key" + [char]0x2019 + " to provide the index of the " + [char]0x2018 + "key" + [char]0x2019 + " in the table. Often hashing algorithms use one-way encryption, which increases the security of the table. Unfortunately, sometimes hashing algorithms can lead to collisions where the hashing algorithm produces the same index value for multiple different keys. These collisions can be handled using a variety of methods, such as rehashing, linear probing, and separate chaining. Separate chaining is the best of these methods, inserting all elements that hash into the same slot index into a linked list. This, however, leads to another list that may have to be search through, and some insecurity around the key."
$insertion2.InsertBefore($bodyText + "`r")
